$d = $word.ActiveDocument

# Locate the paragraph that holds the "公开地址：" label immediately
# followed by the "HYPERLINK https://www.aminer.cn/influencelocality"
# field - this whole paragraph needs to be replaced by three paragraphs
# pointing at the new Quark Drive share link instead.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if (($text -like "*aminer.cn/influencelocality*") -or ($text -like "*公开地址*")) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    Write-Host "Could not find the '公开地址' / aminer link paragraph."
} else {
    $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体"/><w:b/></w:rPr><w:t>公开地址：</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>链接：https://pan.quark.cn/s/37c29307ba85?pwd=Pwmg</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>提取码：Pwmg</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

    [void]$target.Range.InsertXML($xml)
    Write-Host "Updated weibo dataset link paragraph."
}
